$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for affected rows (only columns D, H, J, K, L, M, O, P changed; these
# values reflect the rows being cyclically re-shuffled among each other).
$rowData = @{
    2  = @{ D = 44468; H = "Verde";            J = 500; K = 1800; L = 2000; M = 1920; O = "Provincia de Linares"; P = 1920 }
    3  = @{ D = 44524; H = "Sin especificar";   J = 200; K = 1500; L = 1600; M = 1550; O = "Provincia de Talca";   P = 1550 }
    4  = @{ D = 44519; H = "Sin especificar";   J = 250; K = 1200; L = 1300; M = 1240; O = "Provincia de Linares"; P = 1240 }
    5  = @{ D = 44477; H = "Sin especificar";   J = 500; K = 1400; L = 1500; M = 1460; O = "Provincia de Linares"; P = 1460 }
    7  = @{ D = 44510; H = "Sin especificar";   J = 600; K = 1300; L = 1400; M = 1350; O = "Provincia de Linares"; P = 1350 }
    8  = @{ D = 44526; H = "Sin especificar";   J = 100; K = 1500; L = 1600; M = 1550; O = "Provincia de Linares"; P = 1550 }
    9  = @{ D = 44511; H = "Sin especificar";   J = 600; K = 1300; L = 1400; M = 1350; O = "Provincia de Linares"; P = 1350 }
    10 = @{ D = 44545; H = "Sin especificar";   J = 550; K = 1700; L = 1800; M = 1755; O = "Provincia de Linares"; P = 1755 }
    11 = @{ D = 44489; H = "Sin especificar";   J = 600; K = 1400; L = 1500; M = 1450; O = "Provincia de Linares"; P = 1450 }
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("H$r").Value = $vals.H
    $ws.Range("J$r").Value = $vals.J
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
}
